# "segrigation of qa and stg" — refresh the reroute-request QA test data on
# the Input sheet: the pickUpDate / OrderId / Amount / Tracking# / WayBill
# columns in rows 2, 9 and 11 are updated to a new batch of values.
#
# Values are staged through a scratch cell that is explicitly formatted as
# Text ("@") and then copied with Paste Special > Values into the real
# destination. That keeps every one of these (numeric-looking date /
# order-id / currency / tracking-number) strings stored as literal text,
# exactly like the existing data in these columns, instead of having Excel
# auto-convert them into real dates/numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

$scratch = $ws.Range("BZ1")

function Set-TextValue {
    param($range, $text)
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)
}

# Row 2 (CreateAccount-style LTL record)
Set-TextValue $ws.Range("B2") "12-27-2021"
Set-TextValue $ws.Range("V2") "51525379"
Set-TextValue $ws.Range("X2") "$473.82"
Set-TextValue $ws.Range("Y2") "FCBTXA1000392"
Set-TextValue $ws.Range("Z2") "FCBTXA1000392"

# Row 9
Set-TextValue $ws.Range("B9") "12-27-2021"
Set-TextValue $ws.Range("V9") "51525380"
Set-TextValue $ws.Range("X9") "$66.05"
Set-TextValue $ws.Range("Y9") "1Z44R7R60392001648"
Set-TextValue $ws.Range("Z9") "FCUPSG1011750"

# Row 11
Set-TextValue $ws.Range("B11") "12-27-2021"
Set-TextValue $ws.Range("V11") "51525381"
Set-TextValue $ws.Range("X11") "$180.13"
Set-TextValue $ws.Range("Y11") "1Z44R7R60398241299"
Set-TextValue $ws.Range("Z11") "FCUPSG1011751"

# Clean up the scratch cell so it doesn't linger as visible data / doesn't
# grow the sheet's used range.
$scratch.Clear()

$excel.CutCopyMode = $false
